$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the report title in A1: generation date 2025-10-20 -> 2025-10-21
$ws.Range("A1").Value2 = "萊爾富 工作統計表  篩選月份：202510   (  製表日期:2025-10-21  )"

# 2. Append a new data row (row 79) below the existing last row (78).
#    Row 78 is an "odd" (non zebra-filled) row, so clone the formatting from
#    row 77 instead, which already carries the alternating-row fill/border/
#    wrap combination that the new row 79 uses.
$ws.Range("A77:AK77").Copy()
$ws.Range("A79:AK79").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3. Fill in the values for the new row 79
$ws.Range("A79").Value2 = 77
$ws.Range("B79").Value2 = "維修"
$ws.Range("C79").Value2 = 2025102391

# D79 must stay a *text* value (it's a long numeric case id, stored as a
# string in the source data) - force the Text number format before writing
# it so Excel doesn't silently coerce the digit string into a number.
$ws.Range("D79").NumberFormat = "@"
$ws.Range("D79").Value2 = "13840114102001"

$ws.Range("E79").Value2 = "一般件"
$ws.Range("F79").Value2 = 3840
$ws.Range("G79").Value2 = "北縣天龍店"
$ws.Range("H79").Value2 = "新北市三重區"
$ws.Range("I79").Value2 = "2025-10-20 10:33:24"
$ws.Range("J79").Value2 = "星期一"
$ws.Range("K79").Value2 = "上午"
$ws.Range("L79").Value2 = "HL23"
$ws.Range("M79").Value2 = "HL-TM主機"
$ws.Range("N79").Value2 = 2304
$ws.Range("O79").Value2 = "錢匣損壞"
$ws.Range("P79").Value2 = "門市反應TM1收銀機(TCX800)(抽屜顏色:白色、鑰匙孔位子(右)、鎖頭編號:5001)夾紙鈔彈簧斷掉1個...請台芝到店協助"
$ws.Range("Q79").Value2 = "THILF03840"
$ws.Range("R79").Value2 = "新北一"
$ws.Range("S79").Value2 = "吳宗鴻"
$ws.Range("T79").Value2 = 1
$ws.Range("U79").Value2 = "已完工"
$ws.Range("V79").Value2 = "2025-10-20 10:36:21"
$ws.Range("W79").Value2 = "2025-10-20 14:49:00"
$ws.Range("X79").Value2 = "2025-10-20 15:19:00"
$ws.Range("Y79").Value2 = "2025-10-21 14:36:00"
$ws.Range("Z79").Value2 = 0.5
$ws.Range("AB79").Value2 = "到場處理"
$ws.Range("AC79").Value2 = "更換錢箱`n換下：81Z1000760`n換上：81Z1004553"
$ws.Range("AK79").Value2 = "O"

# Restore D79's formatting (fill/border/alignment) to match the rest of the
# row now that the text value has been committed, so only the cell's value
# type differs from a plain number cell - not its visual style.
$ws.Range("D77").Copy()
$ws.Range("D79").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 4. Turn on word-wrap for P78 and AC78 (their long text now wraps instead
#    of overflowing / getting cut off)
$ws.Range("P78").WrapText = $true
$ws.Range("AC78").WrapText = $true

# 5. Extend the print area to include the newly added row and move the
#    active cell/selection onto it, matching the author's edit.
foreach ($n in $wb.Names) {
    if ($n.Name() -eq "Report!Print_Area") {
        $n.RefersTo = "='Report'!`$A`$1:`$AK`$79"
    }
}
$ws.Range("A79").Select()
